$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'column_name',
    'ticker',
    'fixed_quarter_date',
    'earnings_call_date',
    'Rating',
    'rating_date',
    'Next Rating',
    'Next Rating Date',
    'Previous Rating',
    'Previous Rating Date',
    'next_rating_date_or_end_of_data',
    'credit_rating_year',
    'previous_fixed_quarter_date',
    'days_since_call_on_fixed_quarter',
    'days_since_rating',
    'for_quarter',
    'for_year',
    'transcript',
    'reportedCurrency',
    'acceptedDate_balance_sheet',
    'cashAndCashEquivalents',
    'shortTermInvestments',
    'cashAndShortTermInvestments',
    'netReceivables',
    'inventory_balance_sheet',
    'otherCurrentAssets',
    'totalCurrentAssets',
    'propertyPlantEquipmentNet',
    'goodwill',
    'intangibleAssets',
    'goodwillAndIntangibleAssets',
    'longTermInvestments',
    'taxAssets',
    'otherNonCurrentAssets',
    'totalNonCurrentAssets',
    'otherAssets',
    'totalAssets',
    'accountPayables',
    'shortTermDebt',
    'taxPayables',
    'deferredRevenue',
    'otherCurrentLiabilities',
    'totalCurrentLiabilities',
    'longTermDebt',
    'deferredRevenueNonCurrent',
    'deferredTaxLiabilitiesNonCurrent',
    'otherNonCurrentLiabilities',
    'totalNonCurrentLiabilities',
    'otherLiabilities',
    'capitalLeaseObligations',
    'totalLiabilities',
    'preferredStock',
    'commonStock',
    'retainedEarnings',
    'accumulatedOtherComprehensiveIncomeLoss',
    'othertotalStockholdersEquity',
    'totalStockholdersEquity',
    'totalEquity',
    'totalLiabilitiesAndStockholdersEquity',
    'minorityInterest',
    'totalLiabilitiesAndTotalEquity',
    'totalInvestments',
    'totalDebt',
    'netDebt',
    'acceptedDate_cash_flow_statement',
    'netIncome_cash_flow_statement',
    'depreciationAndAmortization_cash_flow_statement',
    'deferredIncomeTax',
    'stockBasedCompensation',
    'changeInWorkingCapital',
    'accountsReceivables',
    'inventory_cash_flow_statement',
    'accountsPayables',
    'otherWorkingCapital',
    'otherNonCashItems',
    'netCashProvidedByOperatingActivities',
    'investmentsInPropertyPlantAndEquipment',
    'acquisitionsNet',
    'purchasesOfInvestments',
    'salesMaturitiesOfInvestments',
    'otherInvestingActivites',
    'netCashUsedForInvestingActivites',
    'debtRepayment',
    'commonStockIssued',
    'commonStockRepurchased',
    'dividendsPaid',
    'otherFinancingActivites',
    'netCashUsedProvidedByFinancingActivities',
    'effectOfForexChangesOnCash',
    'netChangeInCash',
    'cashAtEndOfPeriod',
    'cashAtBeginningOfPeriod',
    'operatingCashFlow',
    'capitalExpenditure',
    'freeCashFlow',
    'acceptedDate_income_statement',
    'revenue',
    'costOfRevenue',
    'grossProfit',
    'grossProfitRatio',
    'researchAndDevelopmentExpenses',
    'generalAndAdministrativeExpenses',
    'sellingAndMarketingExpenses',
    'sellingGeneralAndAdministrativeExpenses',
    'otherExpenses',
    'operatingExpenses',
    'costAndExpenses',
    'interestIncome',
    'interestExpense',
    'depreciationAndAmortization_income_statement',
    'ebitda',
    'ebitdaratio',
    'operatingIncome',
    'operatingIncomeRatio',
    'totalOtherIncomeExpensesNet',
    'incomeBeforeTax',
    'incomeBeforeTaxRatio',
    'incomeTaxExpense',
    'netIncome_income_statement',
    'netIncomeRatio',
    'eps',
    'epsdiluted',
    'weightedAverageShsOut',
    'weightedAverageShsOutDil',
    'financial_statement_date',
    'marketCap',
    'EBIT',
    'common_plus_preferred_stock',
    'workingCapital',
    'Ratio_A',
    'Ratio_B',
    'Ratio_C',
    'Ratio_D',
    'Ratio_E',
    'Altman_Z',
    'filingDate',
    'currentRatio',
    'quickRatio',
    'cashRatio',
    'returnOnAssets',
    'returnOnEquity',
    'returnOnCapitalEmployed',
    'EBITtoRevenue',
    'debtRatio',
    'debtRatioAlt',
    'debtEquityRatio',
    'equityMultiplier',
    'enterpriseValueMultiplier',
    'operatingCashFlowPerShare',
    'freeCashFlowPerShare',
    'cashPerShare',
    'operatingCashFlowToSales',
    'freeCashFlowToOperatingCashFlow',
    'Altman_Z_diff',
    'Ratio_A_diff',
    'Ratio_B_diff',
    'Ratio_C_diff',
    'Ratio_D_diff',
    'Ratio_E_diff',
    'grossProfitRatio_diff',
    'ebitdaratio_diff',
    'operatingIncomeRatio_diff',
    'incomeBeforeTaxRatio_diff',
    'netIncomeRatio_diff',
    'rating_on_previous_fixed_quarter_date',
    'Investment_Grade',
    'Change Direction Since Last Fixed Quarter Date',
    'Change Since Last Fixed Quarter Date',
    'Sector',
    'currentRatio_diff',
    'quickRatio_diff',
    'cashRatio_diff',
    'returnOnAssets_diff',
    'returnOnEquity_diff',
    'returnOnCapitalEmployed_diff',
    'EBITtoRevenue_diff',
    'debtRatio_diff',
    'debtRatioAlt_diff',
    'debtEquityRatio_diff',
    'equityMultiplier_diff',
    'enterpriseValueMultiplier_diff',
    'operatingCashFlowPerShare_diff',
    'freeCashFlowPerShare_diff',
    'cashPerShare_diff',
    'operatingCashFlowToSales_diff',
    'freeCashFlowToOperatingCashFlow_diff',
    'pos_score_finbert',
    'num_transparency',
    'gf_score',
    'word_count',
    'num_questions',
    'Positiv',
    'Negativ',
    'Strong',
    'Weak',
    'Active',
    'Passive',
    'Ovrst',
    'Undrst',
    'PN',
    'SW',
    'AP',
    'OU',
    'tone',
    'num_q_by_len'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
